# Trade #70 closed at 2026-02-18 00:27:13 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook to reflect:
#  - Trade #98 (HighProbConvergence, row 99 in "All Trades" / row 11 in
#    "HighProbConvergence") closing out with an early exit.
#  - A brand new open Trade #127 (MarketMaking) appended to "All Trades"
#    (row 128) and to the "MarketMaking" strategy sheet (row 48).
#  - Roll-up metrics on "Summary" and "Strategy Status" refreshed to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - top level roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.2
$summary.Range("B4").Value = 0.31
$summary.Range("B5").Value = 0.06
$summary.Range("B6").Value = 98
$summary.Range("B7").Value = 46
$summary.Range("B9").Value = 46.94

# ---------------------------------------------------------------------
# Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.32
$status.Range("D3").Value = 10
$status.Range("E3").Value = 0.33
$status.Range("F3").Value = 0.32
$status.Range("G3").Value = 80

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #98 (row 99) + append trade #127
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(99, 7).Value = 0.67        # G99 Exit Price
$allTrades.Cells.Item(99, 8).Value = "CLOSED"    # H99 Status
$allTrades.Cells.Item(99, 9).Value = 9.8361      # I99 P&L %
$allTrades.Cells.Item(99, 10).Value = 0.06       # J99 P&L $
$allTrades.Cells.Item(99, 11).Value = 100.32     # K99 Capital After
$allTrades.Cells.Item(99, 12).Value = "early_exit" # L99 Exit Reason
$allTrades.Cells.Item(99, 13).Value = 0.13       # M99 Duration (min)

$newTradeRow = 128
$allTrades.Cells.Item($newTradeRow, 1).Value = 127
$c = $allTrades.Cells.Item($newTradeRow, 2)
$c.NumberFormat = "@"
$c.Value = "2026-02-18"
$c.Style = "Normal"
$c = $allTrades.Cells.Item($newTradeRow, 3)
$c.NumberFormat = "@"
$c.Value = "00:27:07"
$c.Style = "Normal"
$allTrades.Cells.Item($newTradeRow, 4).Value = "MarketMaking"
$allTrades.Cells.Item($newTradeRow, 5).Value = "DOWN"
$allTrades.Cells.Item($newTradeRow, 6).Value = 0.61
# G128 (Exit Price) intentionally left blank - trade is still OPEN
$allTrades.Cells.Item($newTradeRow, 8).Value = "OPEN"
$allTrades.Cells.Item($newTradeRow, 9).Value = 0
$allTrades.Cells.Item($newTradeRow, 10).Value = 0
$allTrades.Cells.Item($newTradeRow, 11).Value = 99.40967800952272
# L128 (Exit Reason) intentionally left blank - trade is still OPEN
$allTrades.Cells.Item($newTradeRow, 13).Value = 0
$allTrades.Cells.Item($newTradeRow, 14).Value = 0
$allTrades.Cells.Item($newTradeRow, 15).Value = 0
$allTrades.Cells.Item($newTradeRow, 16).Value = 0.65
$allTrades.Cells.Item($newTradeRow, 17).Value = "Wide spread capture: 392 bps vs avg 280 bps"

# ---------------------------------------------------------------------
# HighProbConvergence sheet - close out the same trade (row 11)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(11, 7).Value = 0.67          # G11 Exit Price
$hpc.Cells.Item(11, 8).Value = "CLOSED"      # H11 Status
$hpc.Cells.Item(11, 9).Value = 9.8361        # I11 P&L %
$hpc.Cells.Item(11, 10).Value = 0.06         # J11 P&L $
$hpc.Cells.Item(11, 11).Value = 100.32       # K11 Capital After
$hpc.Cells.Item(11, 16).Value = "early_exit" # P11 Exit Reason
$hpc.Cells.Item(11, 17).Value = 0.13         # Q11 Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - append the same new trade #127 (row 48)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$newRow = 48
$mm.Cells.Item($newRow, 1).Value = 127
$c = $mm.Cells.Item($newRow, 2)
$c.NumberFormat = "@"
$c.Value = "2026-02-18"
$c.Style = "Normal"
$c = $mm.Cells.Item($newRow, 3)
$c.NumberFormat = "@"
$c.Value = "00:27:07"
$c.Style = "Normal"
$mm.Cells.Item($newRow, 4).Value = "MarketMaking"
$mm.Cells.Item($newRow, 5).Value = "DOWN"
$mm.Cells.Item($newRow, 6).Value = 0.61
# G48 (Exit Price) intentionally left blank - trade is still OPEN
$mm.Cells.Item($newRow, 8).Value = "OPEN"
$mm.Cells.Item($newRow, 9).Value = 0
$mm.Cells.Item($newRow, 10).Value = 0
$mm.Cells.Item($newRow, 11).Value = 99.40967800952272
$mm.Cells.Item($newRow, 12).Value = 0
$mm.Cells.Item($newRow, 13).Value = 0
$mm.Cells.Item($newRow, 14).Value = 0.65
$mm.Cells.Item($newRow, 15).Value = "Wide spread capture: 392 bps vs avg 280 bps"
# P48 (Exit Reason) intentionally left blank - trade is still OPEN
$mm.Cells.Item($newRow, 17).Value = 0
